# Apply the "added Colorado / run-date" update to the About sheet of the
# Maximum Capacity Factor workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# New label identifying the state this copy of the workbook was built for.
$ws.Range("B1").Value = "Colorado"

# New date stamp (12/12/2023) recording when this regional copy was produced.
$ws.Range("C1").Value = (Get-Date -Year 2023 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C1").NumberFormat = "m/d/yyyy"

# Reflect the new window geometry recorded for this save.
$excel.Windows.Item(1).Left = 390
$excel.Windows.Item(1).Top = 390
$excel.Windows.Item(1).Width = 17430
$excel.Windows.Item(1).Height = 17160
